# Fixed some logic errors: correct the Date column (12/10 instead of 12/11)
# and correct the "Time Spent" values for each application entry.
# Row order / Name values stay the same (Netflix, Notepad, Zoom, Outlook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct date for all data rows: 12/10/2020 instead of 12/11/2020
$newDate = Get-Date -Year 2020 -Month 12 -Day 10 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("A2").Value = $newDate
$ws.Range("A3").Value = $newDate
$ws.Range("A4").Value = $newDate
$ws.Range("A5").Value = $newDate

# Correct the "Time Spent" values (stored as text in column C)
$ws.Range("C2").Value = "00:28:00"
$ws.Range("C3").Value = "00:10:00"
$ws.Range("C4").Value = "06:21:16"
$ws.Range("C5").Value = "00:03:45"
